# bugfix 21 - When the report contains accents in the text, the pdf generated
# doesn't contain the correct text.
#
# 1) The "Data" sheet header row (A1:L1) was missing the accent on "column" ->
#    "cólumn" for all 12 header cells.
# 2) The "Summary" sheet's B4 label value is updated from "Another" to
#    "Número" to exercise/verify accented-character rendering.

$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item("Data")
$dataSheet.Range("A1").Value = "cólumn 1"
$dataSheet.Range("B1").Value = "cólumn 2"
$dataSheet.Range("C1").Value = "cólumn 3"
$dataSheet.Range("D1").Value = "cólumn 4"
$dataSheet.Range("E1").Value = "cólumn 5"
$dataSheet.Range("F1").Value = "cólumn 6"
$dataSheet.Range("G1").Value = "cólumn 7"
$dataSheet.Range("H1").Value = "cólumn 8"
$dataSheet.Range("I1").Value = "cólumn 9"
$dataSheet.Range("J1").Value = "cólumn 10"
$dataSheet.Range("K1").Value = "cólumn 11"
$dataSheet.Range("L1").Value = "cólumn 12"

$summarySheet = $wb.Worksheets.Item("Summary")
$summarySheet.Range("B4").Value = "Número"
